$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# --- Row 1 (header row) ---
# Col 1: "Nombre " + proofErr-wrapped "csv"  ->  single run "Nombre csv"
$cell = $tbl.Cell(1,1)
$cell.Range.Delete()
$cell.Range.InsertAfter("Nombre csv")

# Col 4: "Score " + proofErr-wrapped "kaggle"  ->  single run "Score kaggle"
$cell = $tbl.Cell(1,4)
$cell.Range.Delete()
$cell.Range.InsertAfter("Score kaggle")

# --- Row 2 (xgboost_imputing_1) ---
# Col 3: collapse the multi-run / proofErr-split sentence into one run
$cell = $tbl.Cell(2,3)
$cell.Range.Delete()
$cell.Range.InsertAfter("Imputación de 0 (conf>=0.95) y 1 (conf>=0.5) en base a AR + modelo xgboost wen")

# --- Row 3 (xgboost_imputing_1_v2, green text) ---
# Col 3: collapse the multi-run / proofErr-split sentence into one run,
# keeping the green font color (00B050)
$cell = $tbl.Cell(3,3)
$cell.Range.Delete()
$cell.Range.InsertAfter("Imputación de 0 (conf=1) y 1 (conf>=0.5) en base a AR + modelo xgboost wen")
$cell.Range.Font.Color = 5287936

# --- Row 4 (previously empty row) ---
# Col 1: xgboost_imputing_1_v2
$cell = $tbl.Cell(4,1)
$cell.Range.InsertAfter("xgboost_imputing_1_v2")

# Col 2: para cuando no lleguemos a 3 en un dia
$cell = $tbl.Cell(4,2)
$cell.Range.InsertAfter("para cuando no lleguemos a 3 en un dia")

# Col 3: Imputación solo de 0 con confianza 1 + xgboost wen
$cell = $tbl.Cell(4,3)
$cell.Range.InsertAfter("Imputación solo de 0 con confianza 1 + xgboost wen")
